# Add a new "Save" column (H) to the s_vals sheet, matching the format
# of the existing header cells (e.g. column G) and filling in the value
# for the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format of the last existing header cell (G1) onto the
# new header cell (H1) before setting its value.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1

$excel.CutCopyMode = 0
